# Insert a new weekly price row at row 101 (Puerro, Vega Central Mapocho de
# Santiago). Excel's Rows.Insert shifts the existing rows 101-117 down to
# 102-118 and extends the used range / dimension automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(101).Insert()

$ws.Cells.Item(101, 1).Value  = 9
$ws.Cells.Item(101, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(101, 3).Value  = "Metropolitana"
$ws.Cells.Item(101, 4).Value  = 45015
$ws.Cells.Item(101, 5).Value  = 13
$ws.Cells.Item(101, 6).Value  = 100112005
$ws.Cells.Item(101, 7).Value  = "Puerro"
$ws.Cells.Item(101, 8).Value  = "Sin especificar"
$ws.Cells.Item(101, 9).Value  = "Primera"
$ws.Cells.Item(101, 10).Value = 70
$ws.Cells.Item(101, 11).Value = 7000
$ws.Cells.Item(101, 12).Value = 8000
$ws.Cells.Item(101, 13).Value = 7500
$ws.Cells.Item(101, 14).Value = "`$/paquete 20 unidades"
$ws.Cells.Item(101, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(101, 16).Value = 375
$ws.Cells.Item(101, 17).Value = 20
$ws.Cells.Item(101, 18).Value = "Hortaliza"
